# Apply the "Ver 2.22 -> Ver 2.23" / date 2024/09/10 -> 2024/09/11 update.
#
# The commit bumped the deck's visible version/date stamp by one day:
#   - the footer text box on slide 1 ("Ver 2.22 Last updated on 2024/09/10")
#   - the Notes & Handouts master's automatic date footer (2024/9/10)
#
# (Internal bookkeeping such as ppt/revisionInfo.xml's save-counter/timestamp
# and the renumbering of the customXml SharePoint metadata parts are purely
# artifacts that PowerPoint itself regenerates on every save/sync cycle; they
# are not reachable through the Presentation/Shape/TextRange object model and
# are left untouched here.)

$p = $ppt.ActivePresentation

# --- Slide 1: update the "Ver x.xx Last updated on yyyy/mm/dd" text box ---
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "Ver 2.22 Last updated on*") {
            $tr.Text = "Ver 2.23 Last updated on 2024/09/11"
        }
    }
}

# --- Notes Master: bump the automatically-updating date footer ---
$nm = $p.NotesMaster
$hf = $nm.HeadersFooters
$dt = $hf.DateAndTime
$dt.Text = "2024/9/11"
